# Insert a new data row at row 296 (pushing existing rows 296-352 down to 297-353)
# and populate it with the new observation, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 296..352 down by one to make room for the new row.
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new data record.
$ws.Cells.Item(296, 1).Value = 9
$ws.Cells.Item(296, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(296, 3).Value = "Metropolitana"
$ws.Cells.Item(296, 4).Value = 44644
$ws.Cells.Item(296, 5).Value = 13
$ws.Cells.Item(296, 6).Value = 100112032
$ws.Cells.Item(296, 7).Value = "Zapallo italiano"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Primera"
$ws.Cells.Item(296, 10).Value = 70
$ws.Cells.Item(296, 11).Value = 8000
$ws.Cells.Item(296, 12).Value = 8000
$ws.Cells.Item(296, 13).Value = 8000
$ws.Cells.Item(296, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(296, 15).Value = "Provincia de Cachapoal"
$ws.Cells.Item(296, 16).Value = 133
$ws.Cells.Item(296, 17).Value = 60
$ws.Cells.Item(296, 18).Value = "Hortaliza"
